$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2, 3, 4) cycle: row2 -> row4, row3 -> row2, row4 -> row3
# i.e. each row's original content moves "up" one slot (row3 into row2,
# row4 into row3, row2 into row4), as if the first record fell to the bottom.

# Capture current ("before") values of the columns that participate in the rotation.
$cols = @("A","B","D","E","F","G","H","Q","R","Z","AB","AO")

$before2 = @{}
$before3 = @{}
$before4 = @{}
foreach ($col in $cols) {
    $before2[$col] = $ws.Range("${col}2").Value2
    $before3[$col] = $ws.Range("${col}3").Value2
    $before4[$col] = $ws.Range("${col}4").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $before3[$col]
    $ws.Range("${col}3").Value = $before4[$col]
    $ws.Range("${col}4").Value = $before2[$col]
}
